$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 19-21 (data no longer present after the new pomodoro split)
$ws.Rows("19:21").Delete()

# Column H (Break Time) switches from text "HH:MM:SS" to a real date/time
# serial number, formatted the same way as column B (numFmt "YYYY-MM-DD HH:MM:SS").
$ws.Range("H2:H18").NumberFormat = $ws.Range("B2").NumberFormat

$ws.Cells.Item(2, 1).Value2 = "Introduction"
$ws.Cells.Item(2, 2).Value2 = 44572.54166666666
$ws.Cells.Item(2, 3).Value2 = 9
$ws.Cells.Item(2, 4).Value2 = 9
$ws.Cells.Item(2, 5).Value2 = 1
$ws.Cells.Item(2, 6).Value2 = "13:00:00"
$ws.Cells.Item(2, 7).Value2 = "13:45:00"
$ws.Cells.Item(2, 8).Value2 = 44572.58333333334

$ws.Cells.Item(3, 1).Value2 = "Kubernetes Overview"
$ws.Cells.Item(3, 2).Value2 = 44572.54166666666
$ws.Cells.Item(3, 3).Value2 = 22
$ws.Cells.Item(3, 4).Value2 = 31
$ws.Cells.Item(3, 5).Value2 = 1
$ws.Cells.Item(3, 6).Value2 = "13:00:00"
$ws.Cells.Item(3, 7).Value2 = "13:45:00"
$ws.Cells.Item(3, 8).Value2 = 44572.58333333334

$ws.Cells.Item(4, 1).Value2 = "Setup Kubernetes"
$ws.Cells.Item(4, 2).Value2 = 44572.54166666666
$ws.Cells.Item(4, 3).Value2 = 14
$ws.Cells.Item(4, 4).Value2 = 45
$ws.Cells.Item(4, 5).Value2 = 1
$ws.Cells.Item(4, 6).Value2 = "13:00:00"
$ws.Cells.Item(4, 7).Value2 = "13:45:00"
$ws.Cells.Item(4, 8).Value2 = 44572.58333333334

$ws.Cells.Item(5, 1).Value2 = "Setup Kubernetes"
$ws.Cells.Item(5, 2).Value2 = 44572.54166666666
$ws.Cells.Item(5, 3).Value2 = 6
$ws.Cells.Item(5, 4).Value2 = 51
$ws.Cells.Item(5, 5).Value2 = 2
$ws.Cells.Item(5, 6).Value2 = "14:00:00"
$ws.Cells.Item(5, 7).Value2 = "14:45:00"
$ws.Cells.Item(5, 8).Value2 = 44572.625

$ws.Cells.Item(6, 1).Value2 = "Kubernetes Concepts"
$ws.Cells.Item(6, 2).Value2 = 44572.54166666666
$ws.Cells.Item(6, 3).Value2 = 13
$ws.Cells.Item(6, 4).Value2 = 64
$ws.Cells.Item(6, 5).Value2 = 2
$ws.Cells.Item(6, 6).Value2 = "14:00:00"
$ws.Cells.Item(6, 7).Value2 = "14:45:00"
$ws.Cells.Item(6, 8).Value2 = 44572.625

$ws.Cells.Item(7, 1).Value2 = "YAML Introduction"
$ws.Cells.Item(7, 2).Value2 = 44572.54166666666
$ws.Cells.Item(7, 3).Value2 = 8
$ws.Cells.Item(7, 4).Value2 = 72
$ws.Cells.Item(7, 5).Value2 = 2
$ws.Cells.Item(7, 6).Value2 = "14:00:00"
$ws.Cells.Item(7, 7).Value2 = "14:45:00"
$ws.Cells.Item(7, 8).Value2 = 44572.625

$ws.Cells.Item(8, 1).Value2 = "Kubernets Concepts - PODs, ReplicaSets, Deployments"
$ws.Cells.Item(8, 2).Value2 = 44572.54166666666
$ws.Cells.Item(8, 3).Value2 = 18
$ws.Cells.Item(8, 4).Value2 = 90
$ws.Cells.Item(8, 5).Value2 = 2
$ws.Cells.Item(8, 6).Value2 = "14:00:00"
$ws.Cells.Item(8, 7).Value2 = "14:45:00"
$ws.Cells.Item(8, 8).Value2 = 44572.625

$ws.Cells.Item(9, 1).Value2 = "Kubernets Concepts - PODs, ReplicaSets, Deployments"
$ws.Cells.Item(9, 2).Value2 = 44572.54166666666
$ws.Cells.Item(9, 3).Value2 = 45
$ws.Cells.Item(9, 4).Value2 = 135
$ws.Cells.Item(9, 5).Value2 = 3
$ws.Cells.Item(9, 6).Value2 = "15:00:00"
$ws.Cells.Item(9, 7).Value2 = "15:45:00"
$ws.Cells.Item(9, 8).Value2 = 44572.66666666666

$ws.Cells.Item(10, 1).Value2 = "Kubernets Concepts - PODs, ReplicaSets, Deployments"
$ws.Cells.Item(10, 2).Value2 = 44572.54166666666
$ws.Cells.Item(10, 3).Value2 = 45
$ws.Cells.Item(10, 4).Value2 = 180
$ws.Cells.Item(10, 5).Value2 = 4
$ws.Cells.Item(10, 6).Value2 = "16:00:00"
$ws.Cells.Item(10, 7).Value2 = "16:45:00"
$ws.Cells.Item(10, 8).Value2 = 44572.70833333334

$ws.Cells.Item(11, 1).Value2 = "Kubernets Concepts - PODs, ReplicaSets, Deployments"
$ws.Cells.Item(11, 2).Value2 = 44573.54166666666
$ws.Cells.Item(11, 3).Value2 = 20
$ws.Cells.Item(11, 4).Value2 = 200
$ws.Cells.Item(11, 5).Value2 = 5
$ws.Cells.Item(11, 6).Value2 = "13:00:00"
$ws.Cells.Item(11, 7).Value2 = "13:45:00"
$ws.Cells.Item(11, 8).Value2 = 44573.58333333334

$ws.Cells.Item(12, 1).Value2 = "Networking in Kubernetes"
$ws.Cells.Item(12, 2).Value2 = 44573.54166666666
$ws.Cells.Item(12, 3).Value2 = 5
$ws.Cells.Item(12, 4).Value2 = 205
$ws.Cells.Item(12, 5).Value2 = 5
$ws.Cells.Item(12, 6).Value2 = "13:00:00"
$ws.Cells.Item(12, 7).Value2 = "13:45:00"
$ws.Cells.Item(12, 8).Value2 = 44573.58333333334

$ws.Cells.Item(13, 1).Value2 = "Services"
$ws.Cells.Item(13, 2).Value2 = 44573.54166666666
$ws.Cells.Item(13, 3).Value2 = 20
$ws.Cells.Item(13, 4).Value2 = 225
$ws.Cells.Item(13, 5).Value2 = 5
$ws.Cells.Item(13, 6).Value2 = "13:00:00"
$ws.Cells.Item(13, 7).Value2 = "13:45:00"
$ws.Cells.Item(13, 8).Value2 = 44573.58333333334

$ws.Cells.Item(14, 1).Value2 = "Services"
$ws.Cells.Item(14, 2).Value2 = 44573.54166666666
$ws.Cells.Item(14, 3).Value2 = 4
$ws.Cells.Item(14, 4).Value2 = 229
$ws.Cells.Item(14, 5).Value2 = 6
$ws.Cells.Item(14, 6).Value2 = "14:00:00"
$ws.Cells.Item(14, 7).Value2 = "14:45:00"
$ws.Cells.Item(14, 8).Value2 = 44573.625

$ws.Cells.Item(15, 1).Value2 = "Microservices Architechture"
$ws.Cells.Item(15, 2).Value2 = 44573.54166666666
$ws.Cells.Item(15, 3).Value2 = 41
$ws.Cells.Item(15, 4).Value2 = 270
$ws.Cells.Item(15, 5).Value2 = 6
$ws.Cells.Item(15, 6).Value2 = "14:00:00"
$ws.Cells.Item(15, 7).Value2 = "14:45:00"
$ws.Cells.Item(15, 8).Value2 = 44573.625

$ws.Cells.Item(16, 1).Value2 = "Microservices Architechture"
$ws.Cells.Item(16, 2).Value2 = 44573.54166666666
$ws.Cells.Item(16, 3).Value2 = 4
$ws.Cells.Item(16, 4).Value2 = 274
$ws.Cells.Item(16, 5).Value2 = 7
$ws.Cells.Item(16, 6).Value2 = "15:00:00"
$ws.Cells.Item(16, 7).Value2 = "15:45:00"
$ws.Cells.Item(16, 8).Value2 = 44573.66666666666

$ws.Cells.Item(17, 1).Value2 = "Kubernetes on the Cloud"
$ws.Cells.Item(17, 2).Value2 = 44573.54166666666
$ws.Cells.Item(17, 3).Value2 = 26
$ws.Cells.Item(17, 4).Value2 = 300
$ws.Cells.Item(17, 5).Value2 = 7
$ws.Cells.Item(17, 6).Value2 = "15:00:00"
$ws.Cells.Item(17, 7).Value2 = "15:45:00"
$ws.Cells.Item(17, 8).Value2 = 44573.66666666666

$ws.Cells.Item(18, 1).Value2 = "Conclusion"
$ws.Cells.Item(18, 2).Value2 = 44573.54166666666
$ws.Cells.Item(18, 3).Value2 = 2
$ws.Cells.Item(18, 4).Value2 = 302
$ws.Cells.Item(18, 5).Value2 = 7
$ws.Cells.Item(18, 6).Value2 = "15:00:00"
$ws.Cells.Item(18, 7).Value2 = "15:45:00"
$ws.Cells.Item(18, 8).Value2 = 44573.66666666666

